$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Nutzwertanalyse section (rows 28-35) ---
# Two section-header rows are removed ("1. Erfuellung der Anforderungen" in row 28
# and "2. Termineinhaltung" in row 35). Before removing them, drop the left indent
# that used to visually nest the criteria rows (29-34) under the "1. ..." header -
# since that header disappears, the indent formatting on those rows is cleared too.
$ws.Range("A29:A34").IndentLevel = 0

# Clear the stray numeric value that lived in the (about to be removed) header row.
$ws.Range("B28").ClearContents()

# Remove the two header rows - row 35 first so row 28's index stays valid.
$ws.Rows.Item(35).Delete()
$ws.Rows.Item(28).Delete()

# --- View state ---
# Selection / scroll position after the edit (rows shifted up by the deletions).
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A28:A33").Select()
